$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update C61: "Av. San Juan 3960" -> "San Juan 3960"
$ws.Range("C61").Value = "San Juan 3960"

# --- Update row 85 fields
$ws.Range("C85").Value = "Espinosa 591"
$ws.Range("M85").Value = -58.449
$ws.Range("N85").Value = -34.616077
$ws.Range("O85").Value = "Boedo"
$ws.Range("P85").Value = "Capital Sur"

# --- Add new row 87
# Columns that look numeric/date-like need a leading apostrophe to force
# text storage (matching the source file, which stores them as text),
# followed by resetting the style so no stray number-format style sticks.
$ws.Range("A87").Value = "'-513"
$ws.Range("A87").Style = "Normal"

$ws.Range("B87").Value = "'7/15/2025"
$ws.Range("B87").Style = "Normal"

$ws.Range("C87").Value = "Montes de Oca 1809"

$ws.Range("D87").Value = "'4"
$ws.Range("D87").Style = "Normal"

$ws.Range("E87").Value = "'808240768"
$ws.Range("E87").Style = "Normal"

$ws.Range("F87").Value = "Optical Power"
$ws.Range("G87").Value = "Pendiente"
$ws.Range("H87").Value = "Colocar columna donde esta el monoducto para acceso a edifciio"
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = "Cambio"
$ws.Range("K87").Value = "Sin equipos"
$ws.Range("L87").Value = "Pasante"
$ws.Range("M87").Value = -58.372941
$ws.Range("N87").Value = -34.648341
$ws.Range("O87").Value = "San Telmo"
$ws.Range("P87").Value = "Capital Sur"
